$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# Insert a new row at row 42, pushing the existing rows 42-46 down to 43-47.
$ws.Rows.Item(42).Insert()

# Fill in the new payment_type option: 3 classes per week at $500.
$ws.Cells.Item(42, 1).Value = 'payment_type'
$ws.Cells.Item(42, 2).Value = '3_X_SEMANA'
$ws.Cells.Item(42, 3).Value = 'Mensual 3 x Semana $500'

# Match the author's final view/selection/print-setup state.
$ws.Activate()
$ws.PageSetup.Orientation = 1
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("B42").Select()
